$wb = $excel.ActiveWorkbook

# --- "LIST" sheet: fill in the list of codes in column A, add "RO.FOU" ---
$wsList = $wb.Worksheets.Item("LIST")

$wsList.Range("A3").Value = "AD.SEC.001.FON.01"
$wsList.Range("A4").Value = "AD.SEC.001.FON.03"
$wsList.Range("A5").Value = "RO.ACT"

$wsList.Range("A6").NumberFormat = "@"
$wsList.Range("A6").Value = "RO.FOU"

$wsList.Range("A7").Clear()
$wsList.Range("A7").Value = "AD.SEC.014.FON.01"

$wsList.Range("A8:A10").Clear()

# --- "Feuil1" sheet: insert "RO.FOU" before "AD.SEC.014.FON.01" ---
$wsFeuil1 = $wb.Worksheets.Item("Feuil1")

$wsFeuil1.Range("B13").NumberFormat = "@"
$wsFeuil1.Range("B13").Value = "RO.FOU"

$wsFeuil1.Range("B14").Clear()
$wsFeuil1.Range("B14").Value = "AD.SEC.014.FON.01"

# --- selections: update to reflect the new ranges ---
$wsFeuil1.Range("B9:B14").Select()
$wsList.Range("C9:D9").Select()
